# Updated the pinout with uSD and DCMI interfaces
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- MCU Ethernet: PC1 = MDC --------------------------------------------
$ws.Cells.Item(35, 2).Value = "MDC"
$ws.Cells.Item(35, 3).Value = "ETH"

# --- uSD (SDIO) interface -------------------------------------------------
$ws.Cells.Item(26, 3).Value = "uSD"
$ws.Cells.Item(26, 2).Value = "D4"

$ws.Cells.Item(27, 2).Value = "D5"
$ws.Cells.Item(27, 3).Value = "uSD"

$ws.Cells.Item(40, 2).Value = "D6"
$ws.Cells.Item(40, 3).Value = "uSD"

$ws.Cells.Item(41, 2).Value = "D7"
$ws.Cells.Item(41, 3).Value = "uSD"

$ws.Cells.Item(42, 2).Value = "D0"
$ws.Cells.Item(42, 3).Value = "uSD"

$ws.Cells.Item(43, 2).Value = "D2"
$ws.Cells.Item(43, 3).Value = "uSD"

$ws.Cells.Item(44, 2).Value = "D3"
$ws.Cells.Item(44, 3).Value = "uSD"

$ws.Cells.Item(45, 2).Value = "CK"
$ws.Cells.Item(45, 3).Value = "uSD"

$ws.Cells.Item(52, 2).Value = "CMD"
$ws.Cells.Item(52, 3).Value = "uSD"

# --- Camera (DCMI) interface ----------------------------------------------
$ws.Cells.Item(6, 2).Value = "HSYNC"
$ws.Cells.Item(6, 3).Value = "Camera"

$ws.Cells.Item(8, 2).Value = "PIXCK"
$ws.Cells.Item(8, 3).Value = "Camera"

$ws.Cells.Item(11, 2).Value = "D0"
$ws.Cells.Item(11, 3).Value = "Camera"

$ws.Cells.Item(12, 2).Value = "D1"
$ws.Cells.Item(12, 3).Value = "Camera"

$ws.Cells.Item(24, 2).Value = "D5"
$ws.Cells.Item(24, 3).Value = "Camera"

$ws.Cells.Item(25, 2).Value = "VSYNC"
$ws.Cells.Item(25, 3).Value = "Camera"

$ws.Cells.Item(66, 2).Value = "D2"
$ws.Cells.Item(66, 3).Value = "Camera"

$ws.Cells.Item(67, 2).Value = "D3"
$ws.Cells.Item(67, 3).Value = "Camera"

$ws.Cells.Item(70, 2).Value = "D4"
$ws.Cells.Item(70, 3).Value = "Camera"

$ws.Cells.Item(71, 2).Value = "D6"
$ws.Cells.Item(71, 3).Value = "Camera"

$ws.Cells.Item(72, 2).Value = "D7"
$ws.Cells.Item(72, 3).Value = "Camera"

# --- Row heights: the saved sheet explicitly pins the (previously blank)
#     data rows to 15pt; rows that already had a full PORT/Function/Grouping
#     triple (the pre-existing ETH rows) are left untouched. --------------
$ws.Range("A2:C2").Rows.RowHeight = 15
$ws.Range("A5:C8").Rows.RowHeight = 15
$ws.Range("A10:C37").Rows.RowHeight = 15
$ws.Range("A40:C108").Rows.RowHeight = 15
$ws.Range("A110:C110").Rows.RowHeight = 15
$ws.Range("A113:C113").Rows.RowHeight = 15

# --- Final selection, as left by the editor --------------------------------
$ws.Range("B23:C23").Select() | Out-Null
